$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'260.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.01%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'26.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.36%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.704"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.37%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06219"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.18%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'1.11%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8513"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.67%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9143"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.79%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-0.14%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.04951"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.64%"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'-0.22%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03100"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.74%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09051"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.30%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001539"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.18%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006156"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.41%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005974"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.94%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-0.28%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.172"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.77%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.146"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.01%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D21").Value = "'0.1310"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.99%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.095"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.10%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04239"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.22%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-1.53%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004078"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'4.27%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'0.07%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D40").Value = "'0.03964"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.39%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'-0.36%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.004134"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.08%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002152"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-2.53%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-18.17%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'-2.90%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.08%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D48").Value = "'0.2215"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'63.73%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.08%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.08%"
$ws.Range("E50").Style = "Normal"
